# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.182878228561681, 1.65323645889881,  0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    3 = @(3.182878228561681, 1.65323645889881,  3.082599426703578,  0.4998867070740569, 8.418600821238126)
    4 = @(0.7287194209349384, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 3.034748368925986)
    5 = @(0.7287194209349384, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 3.594575437922795)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]  # B - TB
    $ws.Cells.Item($row, 3).Value = $vals[1]  # C - d2S
    $ws.Cells.Item($row, 4).Value = $vals[2]  # D - K
    $ws.Cells.Item($row, 5).Value = $vals[3]  # E - IP
    $ws.Cells.Item($row, 7).Value = $vals[4]  # G - sum
}
